$d = $word.ActiveDocument

# 1) Replace the placeholder-ID text in the first paragraph. The old content
#    was split across two runs ("**ID__..._topic_25__ID**" + a trailing
#    space run); Find/Execute over the combined range collapses that into a
#    single run carrying the new ID text.
$d.Content.Find.Execute(
    "**ID__AFFARS_pgi_5315_topic_25__ID** ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "**ID__AFFARS_AFMC_PGI_5315_404_1_90__ID**", 2) | Out-Null

# 2) Update the first paragraph's formatting: add a paragraph border (5pt
#    space on every edge) and widen the left indent from 120 twips (6pt) to
#    225 twips (11.25pt).
$p = $d.Paragraphs(1)
$borders = $p.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
$p.Range.ParagraphFormat.LeftIndent = 11.25
